# Update column F (dSF) values to match repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -1
    8  = -4
    9  = 1
    10 = 3
    11 = 2
    12 = 4
    13 = -5
    15 = 2
    16 = 4
    18 = -1
    20 = -5
    21 = 6
    22 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
